# Tambah lap AP1, security dan perbaiki logika logout
# Inserts a new column "KODE REK" before the existing "NAMA PAKET / KEGIATAN"
# column, renames "NO." -> "NO", fixes up the row-number guide row, column
# widths, merges and the active selection to match the edited template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand-new column at B (shifts old B:K -> C:L, formulas/merges
#    shift automatically, same as a user doing Right-click > Insert on the
#    column header).
$ws.Columns.Item(2).Insert()

# 2) Column widths: column A becomes the narrow "NO" column, column B takes
#    over the width the old column A used to have, and column E (used to be
#    D, "SUMBER DANA") gets a bit narrower to make room.
$ws.Columns.Item(1).ColumnWidth = 4
$ws.Columns.Item(2).ColumnWidth = 10.6328125
$ws.Columns.Item(5).ColumnWidth = 15.36328125

# 3) The label rows (3-7) don't use column B at all - clear whatever
#    formatting the insert copied into it so it stays a plain empty cell.
$ws.Range("B3:B7").ClearFormats()

# 4) Header row fix-ups.
$ws.Range("A8").Value = "NO"
$ws.Range("B8").Value = "KODE REK"
$ws.Range("B8").Style = $ws.Range("A8").Style
$ws.Range("B8:B9").Merge()

# 5) Row 10 is the little "1,2,3..." column-reference guide under the
#    headers - renumber it in full now that there are 12 columns.
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 6
$ws.Range("G10").Value = 7
$ws.Range("H10").Value = 8
$ws.Range("I10").Value = 9
$ws.Range("J10").Value = 10
$ws.Range("K10").Value = 11
$ws.Range("L10").Value = 12
$ws.Range("B10").Style = $ws.Range("A10").Style
$ws.Range("E10").Style = $ws.Range("A10").Style
$ws.Range("K10").Style = $ws.Range("F10").Style

# 6) Row 11 (blank footer row) - E11 matches the A11 style, K11 matches F11.
$ws.Range("E11").Style = $ws.Range("A11").Style
$ws.Range("K11").Style = $ws.Range("F11").Style

# 7) Selection moves to A11 and the frozen/top-left cell scroll offset is
#    cleared.
$ws.Range("A11").Select()
